$wb = $excel.ActiveWorkbook

# --- Sheet: ev_charging_uc ---
# C13 / C14 hold long comma-separated timeslice lists; G7 (=C14) and G8 (=C13)
# cache their formula results. The lists were re-shuffled (same tokens, new order).
$wsEv = $wb.Worksheets.Item("ev_charging_uc")
$wsEv.Range("C14").Value = 'S6b1008h04,S6b1010h20,S6b1011h02,S6b1012h19,S6b1013h02,S6c1217h02,S6c1217h06,S6c1217h24,S6c1221h19,S6c1221h21,S6c1222h19,S1aH1,S6b1009h02,S6b1013h19,S6c1217h21,S6c1220h06,S6c1221h22,S2aH1,S5aH8,S6b1009h19,S6b1011h04,S6b1013h05,S6b1013h20,S6b1013h23,S6c1217h04,S6c1220h01,S6c1220h19,S6c1221h24,S6c1222h05,S3aH8,S6b1008h02,S6b1008h21,S6b1009h03,S6b1009h24,S6b1010h21,S6b1010h24,S6b1012h04,S6b1012h21,S6b1013h04,S6c1218h04,S6c1218h06,S6c1221h05,S6c1221h20,S6c1223h03,S6aH8,S6b1008h19,S6b1009h06,S6b1012h22,S6b1013h03,S6b1014h01,S6c1222h02,S6c1222h06,S4aH1,S6aH1,S6b1008h24,S6b1010h01,S6b1011h05,S6b1011h06,S6b1011h22,S6b1014h03,S6c1218h23,S6c1219h04,S6c1223h22,S6b1008h22,S6b1009h21,S6b1009h23,S6b1010h23,S6b1012h03,S6b1014h04,S6b1014h05,S6b1014h19,S6c1217h03,S6c1219h01,S6b1009h04,S6b1010h19,S6b1011h03,S6b1011h19,S6b1014h24,S6c1218h02,S6c1218h21,S6c1218h22,S6c1219h02,S6c1219h06,S6c1220h02,S6b1008h03,S6b1012h23,S6b1013h24,S6b1014h06,S6b1014h21,S6c1218h01,S6c1219h20,S6c1220h22,S6c1223h20,S6b1010h22,S6b1011h24,S6b1012h06,S6b1014h23,S6c1218h20,S6c1218h24,S6c1219h05,S6c1219h23,S6c1220h03,S6c1221h03,S6c1221h23,S6c1222h22,S6c1223h06,S6c1223h19,S6b1010h03,S6b1011h01,S6b1012h24,S6c1217h20,S6c1218h03,S6c1220h24,S6c1221h04,S6c1222h21,S6c1223h02,S6c1223h05,S1aH8,S5aH1,S6b1014h22,S6c1217h19,S6c1219h19,S6c1219h21,S6c1220h21,S6c1221h01,S6c1221h02,S6c1222h03,S4aH8,S6b1008h05,S6b1008h20,S6b1010h04,S6b1010h06,S6b1011h20,S6b1011h23,S6b1012h02,S6b1013h06,S6b1013h22,S6b1014h02,S6c1217h22,S6c1217h23,S6c1219h03,S6c1219h24,S6c1221h06,S6c1222h01,S6c1222h20,S6c1222h23,S6c1223h04,S6c1223h23,S6b1009h05,S6b1010h02,S6b1011h21,S6b1012h01,S6b1013h21,S6b1014h20,S6c1217h01,S6c1220h05,S3aH1,S6b1009h20,S6b1010h05,S6b1012h05,S6b1012h20,S6c1218h05,S6c1218h19,S6c1219h22,S6c1220h04,S6c1220h20,S6c1220h23,S6c1222h04,S6c1222h24,S6c1223h21,S2aH8,S6b1008h01,S6b1008h06,S6b1008h23,S6b1009h01,S6b1009h22,S6b1013h01,S6c1217h05,S6c1223h01,S6c1223h24'
$wsEv.Range("C13").Value = 'S1aH3,S6b1008h16,S6b1010h09,S6b1011h14,S6b1012h14,S6b1012h15,S6c1217h10,S6c1217h18,S6c1220h16,S6c1221h09,S6c1222h07,S6b1012h12,S6b1013h15,S6b1013h16,S6b1014h18,S6c1217h07,S6c1218h12,S6c1218h17,S6c1219h11,S6c1219h18,S6c1220h10,S6c1221h07,S6c1222h09,S6c1222h12,S6c1223h09,S4aH2,S6b1010h15,S6b1012h07,S6b1013h18,S6b1014h13,S6c1219h16,S6c1221h11,S3aH2,S6b1009h17,S6b1011h07,S6b1011h11,S6b1013h10,S6c1217h16,S6c1219h17,S6c1220h18,S6c1221h13,S1aH6,S6b1008h12,S6b1008h13,S6b1010h07,S6b1010h13,S6b1010h14,S6b1011h12,S6b1012h11,S6b1012h18,S6c1217h09,S6c1218h10,S6c1218h11,S6c1218h16,S6c1219h08,S6c1223h08,S6c1223h16,S1aH4,S2aH2,S4aH6,S6b1008h17,S6b1010h16,S6b1012h13,S6b1013h08,S6b1014h09,S6c1218h07,S6c1220h09,S6c1220h13,S6c1220h15,S6c1221h10,S6c1222h13,S6c1222h17,S5aH4,S6aH3,S6aH4,S6b1008h07,S6b1008h11,S6b1010h18,S6b1013h09,S6b1013h14,S6c1218h15,S6c1219h07,S6c1220h17,S6c1222h16,S6c1223h18,S1aH5,S2aH4,S2aH5,S5aH2,S5aH5,S6aH2,S6b1008h15,S6b1010h08,S6b1011h16,S6b1013h11,S6b1014h15,S6c1223h12,S6c1223h15,S1aH2,S1aH7,S2aH3,S3aH6,S4aH7,S5aH6,S5aH7,S6aH6,S6b1008h08,S6b1009h07,S6b1009h12,S6b1009h14,S6b1011h09,S6b1012h09,S6b1013h12,S6c1217h11,S6c1217h13,S6c1218h13,S6c1221h14,S6c1221h15,S6c1221h17,S6c1221h18,S6c1222h11,S4aH5,S6b1008h09,S6b1011h08,S6b1011h13,S6b1014h07,S6b1014h12,S6b1014h14,S6b1014h16,S6c1220h08,S6c1221h08,S6c1221h16,S6c1222h10,S2aH7,S3aH4,S4aH3,S4aH4,S6b1009h18,S6b1010h10,S6b1011h15,S6b1011h18,S6b1012h16,S6b1012h17,S6c1218h08,S6c1218h18,S6c1219h10,S6c1220h11,S6c1222h08,S6c1223h07,S2aH6,S3aH5,S5aH3,S6b1009h11,S6b1009h15,S6b1014h17,S6c1217h15,S6c1218h14,S6c1219h09,S6c1220h07,S6c1220h12,S6c1222h15,S6c1223h11,S6aH5,S6b1011h10,S6b1014h11,S6c1217h08,S6c1217h14,S6c1219h12,S6c1219h15,S6c1221h12,S6c1222h14,S6c1223h10,S6c1223h13,S6c1223h14,S6c1223h17,S3aH3,S6b1010h17,S6b1011h17,S6b1012h10,S6b1013h13,S6b1013h17,S6c1217h12,S6c1218h09,S6c1219h14,S6aH7,S6b1008h14,S6b1009h08,S6b1009h10,S6b1009h13,S6b1009h16,S6b1010h11,S6b1013h07,S6b1014h08,S6b1014h10,S6c1219h13,S6c1220h14,S6c1222h18,S3aH7,S6b1008h10,S6b1008h18,S6b1009h09,S6b1010h12,S6b1012h08,S6c1217h17'

# --- Sheet: re_profiles ---
$wsRe = $wb.Worksheets.Item("re_profiles")

# Rename shared string "elc_win-ITA" -> "elc_won-ITA" across every cell that uses it
# (column K, rows 4-387) so the shared string table collapses to the new text
# instead of forking a duplicate entry.
$wsRe.Range("K4:K387").Value = "elc_won-ITA"

# Re-sorted random sample rows (label + value move together):
$wsRe.Range("M4").Value = "S4"
$wsRe.Range("N4").Value = 0.25252198056455344

$wsRe.Range("M5").Value = "S1"
$wsRe.Range("N5").Value = 0.14770939379916706

$wsRe.Range("M6").Value = "S3"
$wsRe.Range("N6").Value = 0.37937297547431742

$wsRe.Range("M9").Value = "S6"
$wsRe.Range("N9").Value = 0.25116844053678855
